$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the species text in B3 from "Many species of phytoplankton" to
# "Alexandrium spp. (dinoflagellate)", with "Alexandrium " in italics and
# the remainder in the regular (non-italic) font, matching the styling
# already used for similar species cells (B2, B4).
$ws.Range("B3").Value = "Alexandrium spp. (dinoflagellate)"
$ws.Range("B3").Characters(1, 12).Font.Italic = $true
$ws.Range("B3").Characters(13, 21).Font.Italic = $false
$ws.Range("B3").Characters(13, 21).Font.Size = 12
$ws.Range("B3").Characters(13, 21).Font.Name = "Calibri"

# Widen column B slightly to better fit the new text.
$ws.Columns.Item(2).ColumnWidth = 28.33

# Update the active selection to D5.
$null = $ws.Range("D5").Select()
